$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "model_3_6_24"
$ws.Cells.Item(2, 2).Value = 0.422046706849045
$ws.Cells.Item(2, 3).Value = 0.1635538922638263
$ws.Cells.Item(2, 4).Value = -0.2097555963164301
$ws.Cells.Item(2, 5).Value = 0.1378663516932161
$ws.Cells.Item(2, 6).Value = 0.6396239995956421
$ws.Cells.Item(2, 7).Value = 0.7474936246871948
$ws.Cells.Item(2, 8).Value = 0.6685941219329834
$ws.Cells.Item(2, 9).Value = 0.7103652954101562
$ws.Cells.Item(3, 1).Value = "model_3_6_23"
$ws.Cells.Item(3, 2).Value = 0.4221359513647231
$ws.Cells.Item(3, 3).Value = 0.1640866229152774
$ws.Cells.Item(3, 4).Value = -0.2093874066791481
$ws.Cells.Item(3, 5).Value = 0.1382905114328093
$ws.Cells.Item(3, 6).Value = 0.639525294303894
$ws.Cells.Item(3, 7).Value = 0.7470175623893738
$ws.Cells.Item(3, 8).Value = 0.6683906316757202
$ws.Cells.Item(3, 9).Value = 0.7100158929824829
$ws.Cells.Item(4, 1).Value = "model_3_6_22"
$ws.Cells.Item(4, 2).Value = 0.4225156818447938
$ws.Cells.Item(4, 3).Value = 0.1662522767638308
$ws.Cells.Item(4, 4).Value = -0.2078421606982517
$ws.Cells.Item(4, 5).Value = 0.1400202496746791
$ws.Cells.Item(4, 6).Value = 0.639104962348938
$ws.Cells.Item(4, 7).Value = 0.7450821995735168
$ws.Cells.Item(4, 8).Value = 0.6675366163253784
$ws.Cells.Item(4, 9).Value = 0.7085906267166138
$ws.Cells.Item(5, 1).Value = "model_3_6_21"
$ws.Cells.Item(5, 2).Value = 0.4227828186365898
$ws.Cells.Item(5, 3).Value = 0.1682271962213757
$ws.Cells.Item(5, 4).Value = -0.2076232231181669
$ws.Cells.Item(5, 5).Value = 0.1412230966446606
$ws.Cells.Item(5, 6).Value = 0.6388093233108521
$ws.Cells.Item(5, 7).Value = 0.7433173656463623
$ws.Cells.Item(5, 8).Value = 0.6674156188964844
$ws.Cells.Item(5, 9).Value = 0.7075994610786438
$ws.Cells.Item(6, 1).Value = "model_3_6_19"
$ws.Cells.Item(6, 2).Value = 0.4228176335876341
$ws.Cells.Item(6, 3).Value = 0.1691297988083424
$ws.Cells.Item(6, 4).Value = -0.2083506012861234
$ws.Cells.Item(6, 5).Value = 0.1415129794264643
$ws.Cells.Item(6, 6).Value = 0.6387708783149719
$ws.Cells.Item(6, 7).Value = 0.7425106763839722
$ws.Cells.Item(6, 8).Value = 0.6678175926208496
$ws.Cells.Item(6, 9).Value = 0.7073606252670288
$ws.Cells.Item(7, 1).Value = "model_3_6_20"
$ws.Cells.Item(7, 2).Value = 0.4230415068136891
$ws.Cells.Item(7, 3).Value = 0.1696466619763916
$ws.Cells.Item(7, 4).Value = -0.2063297103796415
$ws.Cells.Item(7, 5).Value = 0.1424460963894006
$ws.Cells.Item(7, 6).Value = 0.6385231018066406
$ws.Cells.Item(7, 7).Value = 0.7420488595962524
$ws.Cells.Item(7, 8).Value = 0.6667007207870483
$ws.Cells.Item(7, 9).Value = 0.706591784954071
$ws.Cells.Item(8, 1).Value = "model_3_6_18"
$ws.Cells.Item(8, 2).Value = 0.4271133273920696
$ws.Cells.Item(8, 3).Value = 0.1674020637554847
$ws.Cells.Item(8, 4).Value = -0.1446426572292241
$ws.Cells.Item(8, 5).Value = 0.1606305617464169
$ws.Cells.Item(8, 6).Value = 0.6340166926383972
$ws.Cells.Item(8, 7).Value = 0.7440546751022339
$ws.Cells.Item(8, 8).Value = 0.6326082348823547
$ws.Cells.Item(8, 9).Value = 0.6916084885597229
$ws.Cells.Item(9, 1).Value = "model_3_6_17"
$ws.Cells.Item(9, 2).Value = 0.4274805211223944
$ws.Cells.Item(9, 3).Value = 0.1689349453649821
$ws.Cells.Item(9, 4).Value = -0.1421389512445563
$ws.Cells.Item(9, 5).Value = 0.1622994193338129
$ws.Cells.Item(9, 6).Value = 0.6336103081703186
$ws.Cells.Item(9, 7).Value = 0.7426849007606506
$ws.Cells.Item(9, 8).Value = 0.631224513053894
$ws.Cells.Item(9, 9).Value = 0.6902334690093994
$ws.Cells.Item(10, 1).Value = "model_3_6_16"
$ws.Cells.Item(10, 2).Value = 0.4281249376881703
$ws.Cells.Item(10, 3).Value = 0.1715780037277636
$ws.Cells.Item(10, 4).Value = -0.1378158972121635
$ws.Cells.Item(10, 5).Value = 0.1651830514894939
$ws.Cells.Item(10, 6).Value = 0.6328971982002258
$ws.Cells.Item(10, 7).Value = 0.7403228878974915
$ws.Cells.Item(10, 8).Value = 0.6288352608680725
$ws.Cells.Item(10, 9).Value = 0.6878573298454285
$ws.Cells.Item(11, 1).Value = "model_3_6_13"
$ws.Cells.Item(11, 2).Value = 0.4283166729595377
$ws.Cells.Item(11, 3).Value = 0.1621359493033299
$ws.Cells.Item(11, 4).Value = -0.1173115827892444
$ws.Cells.Item(11, 5).Value = 0.1662321233422024
$ws.Cells.Item(11, 6).Value = 0.6326850056648254
$ws.Cells.Item(11, 7).Value = 0.7487608194351196
$ws.Cells.Item(11, 8).Value = 0.6175032258033752
$ws.Cells.Item(11, 9).Value = 0.6869930028915405
$ws.Cells.Item(12, 1).Value = "model_3_6_14"
$ws.Cells.Item(12, 2).Value = 0.4283738870494492
$ws.Cells.Item(12, 3).Value = 0.1637943374490182
$ws.Cells.Item(12, 4).Value = -0.1196467787604738
$ws.Cells.Item(12, 5).Value = 0.1664467540520836
$ws.Cells.Item(12, 6).Value = 0.6326216459274292
$ws.Cells.Item(12, 7).Value = 0.7472787499427795
$ws.Cells.Item(12, 8).Value = 0.618793785572052
$ws.Cells.Item(12, 9).Value = 0.6868160963058472
$ws.Cells.Item(13, 1).Value = "model_3_6_15"
$ws.Cells.Item(13, 2).Value = 0.4283900938383601
$ws.Cells.Item(13, 3).Value = 0.164253155625449
$ws.Cells.Item(13, 4).Value = -0.1203877604601391
$ws.Cells.Item(13, 5).Value = 0.1664759279294605
$ws.Cells.Item(13, 6).Value = 0.632603645324707
$ws.Cells.Item(13, 7).Value = 0.7468687891960144
$ws.Cells.Item(13, 8).Value = 0.6192033290863037
$ws.Cells.Item(13, 9).Value = 0.6867920756340027
$ws.Cells.Item(14, 1).Value = "model_3_6_12"
$ws.Cells.Item(14, 2).Value = 0.4306801642224786
$ws.Cells.Item(14, 3).Value = 0.166503408615174
$ws.Cells.Item(14, 4).Value = -0.09128867704021415
$ws.Cells.Item(14, 5).Value = 0.1769556786428186
$ws.Cells.Item(14, 6).Value = 0.6300693154335022
$ws.Cells.Item(14, 7).Value = 0.7448577880859375
$ws.Cells.Item(14, 8).Value = 0.6031211614608765
$ws.Cells.Item(14, 9).Value = 0.6781572699546814
$ws.Cells.Item(15, 1).Value = "model_3_6_11"
$ws.Cells.Item(15, 2).Value = 0.4326176379169558
$ws.Cells.Item(15, 3).Value = 0.1651193255928873
$ws.Cells.Item(15, 4).Value = -0.06016240340841073
$ws.Cells.Item(15, 5).Value = 0.185984153333383
$ws.Cells.Item(15, 6).Value = 0.6279250383377075
$ws.Cells.Item(15, 7).Value = 0.7460947036743164
$ws.Cells.Item(15, 8).Value = 0.585918664932251
$ws.Cells.Item(15, 9).Value = 0.6707180738449097
$ws.Cells.Item(16, 1).Value = "model_3_6_10"
$ws.Cells.Item(16, 2).Value = 0.4332713766358854
$ws.Cells.Item(16, 3).Value = 0.1517815274669997
$ws.Cells.Item(16, 4).Value = -0.02304291698233363
$ws.Cells.Item(16, 5).Value = 0.190043419932598
$ws.Cells.Item(16, 6).Value = 0.6272015571594238
$ws.Cells.Item(16, 7).Value = 0.7580140233039856
$ws.Cells.Item(16, 8).Value = 0.5654038190841675
$ws.Cells.Item(16, 9).Value = 0.6673734188079834
$ws.Cells.Item(17, 1).Value = "model_3_6_3"
$ws.Cells.Item(17, 2).Value = 0.4342288060946472
$ws.Cells.Item(17, 3).Value = 0.1378714347270967
$ws.Cells.Item(17, 4).Value = 0.07639048276154536
$ws.Cells.Item(17, 5).Value = 0.2134419786637162
$ws.Cells.Item(17, 6).Value = 0.6261419653892517
$ws.Cells.Item(17, 7).Value = 0.770444929599762
$ws.Cells.Item(17, 8).Value = 0.5104501247406006
$ws.Cells.Item(17, 9).Value = 0.6480938196182251
$ws.Cells.Item(18, 1).Value = "model_3_6_9"
$ws.Cells.Item(18, 2).Value = 0.4354447942072245
$ws.Cells.Item(18, 3).Value = 0.1540196448654777
$ws.Cells.Item(18, 4).Value = 0.005864104001596759
$ws.Cells.Item(18, 5).Value = 0.2004526074830838
$ws.Cells.Item(18, 6).Value = 0.6247963309288025
$ws.Cells.Item(18, 7).Value = 0.7560139298439026
$ws.Cells.Item(18, 8).Value = 0.54942786693573
$ws.Cells.Item(18, 9).Value = 0.6587966680526733
$ws.Cells.Item(19, 1).Value = "model_3_6_7"
$ws.Cells.Item(19, 2).Value = 0.4365649007455983
$ws.Cells.Item(19, 3).Value = 0.139257376302887
$ws.Cells.Item(19, 4).Value = 0.0528218906839184
$ws.Cells.Item(19, 5).Value = 0.2067979641459706
$ws.Cells.Item(19, 6).Value = 0.6235566139221191
$ws.Cells.Item(19, 7).Value = 0.7692062854766846
$ws.Cells.Item(19, 8).Value = 0.5234757661819458
$ws.Cells.Item(19, 9).Value = 0.6535682678222656
$ws.Cells.Item(20, 1).Value = "model_3_6_8"
$ws.Cells.Item(20, 2).Value = 0.4366637409533113
$ws.Cells.Item(20, 3).Value = 0.1456140689688583
$ws.Cells.Item(20, 4).Value = 0.04211529083524534
$ws.Cells.Item(20, 5).Value = 0.2070710652526173
$ws.Cells.Item(20, 6).Value = 0.6234472393989563
$ws.Cells.Item(20, 7).Value = 0.763525664806366
$ws.Cells.Item(20, 8).Value = 0.5293929576873779
$ws.Cells.Item(20, 9).Value = 0.6533432602882385
$ws.Cells.Item(21, 1).Value = "model_3_6_4"
$ws.Cells.Item(21, 2).Value = 0.4382588698781968
$ws.Cells.Item(21, 3).Value = 0.1786569187780632
$ws.Cells.Item(21, 4).Value = 0.04889021098287472
$ws.Cells.Item(21, 5).Value = 0.2281797193118984
$ws.Cells.Item(21, 6).Value = 0.6216818690299988
$ws.Cells.Item(21, 7).Value = 0.7339967489242554
$ws.Cells.Item(21, 8).Value = 0.5256486535072327
$ws.Cells.Item(21, 9).Value = 0.6359505653381348
$ws.Cells.Item(22, 1).Value = "model_3_6_2"
$ws.Cells.Item(22, 2).Value = 0.4389482071799513
$ws.Cells.Item(22, 3).Value = 0.1409397237633024
$ws.Cells.Item(22, 4).Value = 0.1610666648279447
$ws.Cells.Item(22, 5).Value = 0.2419306790413518
$ws.Cells.Item(22, 6).Value = 0.6209190487861633
$ws.Cells.Item(22, 7).Value = 0.7677028775215149
$ws.Cells.Item(22, 8).Value = 0.4636522531509399
$ws.Cells.Item(22, 9).Value = 0.624620258808136
$ws.Cells.Item(23, 1).Value = "model_3_6_0"
$ws.Cells.Item(23, 2).Value = 0.4394307493511261
$ws.Cells.Item(23, 3).Value = 0.1136320513434056
$ws.Cells.Item(23, 4).Value = 0.2704366687692686
$ws.Cells.Item(23, 5).Value = 0.2607748648733945
$ws.Cells.Item(23, 6).Value = 0.6203849911689758
$ws.Cells.Item(23, 7).Value = 0.7921064496040344
$ws.Cells.Item(23, 8).Value = 0.4032068848609924
$ws.Cells.Item(23, 9).Value = 0.6090933680534363
$ws.Cells.Item(24, 1).Value = "model_3_6_6"
$ws.Cells.Item(24, 2).Value = 0.4400703738380892
$ws.Cells.Item(24, 3).Value = 0.1610709755207625
$ws.Cells.Item(24, 4).Value = 0.0665626057095523
$ws.Cells.Item(24, 5).Value = 0.2236594402645716
$ws.Cells.Item(24, 6).Value = 0.6196771264076233
$ws.Cells.Item(24, 7).Value = 0.7497125267982483
$ws.Cells.Item(24, 8).Value = 0.5158816576004028
$ws.Cells.Item(24, 9).Value = 0.6396750807762146
$ws.Cells.Item(25, 1).Value = "model_3_6_5"
$ws.Cells.Item(25, 2).Value = 0.4420657517144443
$ws.Cells.Item(25, 3).Value = 0.1381179576395907
$ws.Cells.Item(25, 4).Value = 0.1448874083812709
$ws.Cells.Item(25, 5).Value = 0.2352043293956448
$ws.Cells.Item(25, 6).Value = 0.6174688339233398
$ws.Cells.Item(25, 7).Value = 0.7702245712280273
$ws.Cells.Item(25, 8).Value = 0.4725940227508545
$ws.Cells.Item(25, 9).Value = 0.6301625370979309
$ws.Cells.Item(26, 1).Value = "model_3_6_1"
$ws.Cells.Item(26, 2).Value = 0.4430552000915807
$ws.Cells.Item(26, 3).Value = 0.1233781379621547
$ws.Cells.Item(26, 4).Value = 0.2885435262536972
$ws.Cells.Item(26, 5).Value = 0.2720852264960085
$ws.Cells.Item(26, 6).Value = 0.6163737177848816
$ws.Cells.Item(26, 7).Value = 0.7833969593048096
$ws.Cells.Item(26, 8).Value = 0.3931997716426849
$ws.Cells.Item(26, 9).Value = 0.5997740626335144
